# Rename the sheets (Sheet1 -> m0, Sheet2 -> m1).
# Renaming updates the "gs_rad_tau0" defined name's sheet-qualified
# reference automatically (Sheet1!... -> m0!...).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Name = "m0"
$ws2.Name = "m1"

# Make m0 the active sheet/tab, with C2:C22 selected (was A7 on Sheet1).
$ws1.Activate() | Out-Null
$ws1.Range("C2:C22").Select() | Out-Null
